# Auto-generated script applying the Sagittarius_Profits market-price refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2232
$ws.Range("I12").Value = 348.5
$ws.Range("K12").Value = 348.5
$ws.Range("M12").Value = -178.5
$ws.Range("H80").Value = 712.875
$ws.Range("J80").Value = 671.8570999999999
$ws.Range("L80").Value = 2015.5713
$ws.Range("N80").Value = -4011.5713
$ws.Range("H83").Value = 712.875
$ws.Range("J83").Value = 671.8570999999999
$ws.Range("L83").Value = 6046.7139
$ws.Range("N83").Value = -16030.7139
$ws.Range("H96").Value = 12469.556
$ws.Range("I96").Value = 15175.143
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 45525.429
$ws.Range("L96").Value = 9000
$ws.Range("M96").Value = -44152.429
$ws.Range("N96").Value = -11746
$ws.Range("H101").Value = 12502720
$ws.Range("I101").Value = 25004574
$ws.Range("J101").Value = 866
$ws.Range("K101").Value = 75013722
$ws.Range("L101").Value = 2598
$ws.Range("M101").Value = -75012100
$ws.Range("N101").Value = -5842
$ws.Range("H106").Value = 1045
$ws.Range("I106").Value = 1000
$ws.Range("K106").Value = 1000
$ws.Range("M106").Value = -369

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4532946
$ws.Range("I32").Value = 4378755
$ws.Range("K32").Value = 4378755
$ws.Range("M32").Value = -4378468
$ws.Range("H45").Value = 7625
$ws.Range("I45").Value = 7625
$ws.Range("K45").Value = 7625
$ws.Range("M45").Value = -7248
$ws.Range("H46").Value = 4829.2
$ws.Range("J46").Value = 4829.2
$ws.Range("L46").Value = 4829.2
$ws.Range("N46").Value = -5467.2
$ws.Range("H110").Value = 2643965.2
$ws.Range("I110").Value = 3364538.5
$ws.Range("J110").Value = 1862.6666
$ws.Range("K110").Value = 3364538.5
$ws.Range("L110").Value = 1862.6666
$ws.Range("M110").Value = -3362493.5
$ws.Range("N110").Value = -5952.6666
$ws.Range("H122").Value = 1100
$ws.Range("I122").Value = 1050
$ws.Range("K122").Value = 3150
$ws.Range("M122").Value = -700
$ws.Range("H132").Value = 1343.75
$ws.Range("I132").Value = 1391.6666
$ws.Range("K132").Value = 4174.9998
$ws.Range("M132").Value = -1644.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J99").Value = 1300
$ws.Range("L99").Value = 1300
$ws.Range("N99").Value = -4296
$ws.Range("H107").Value = 1341.7142
$ws.Range("I107").Value = 1224.75
$ws.Range("J107").Value = 1497.6666
$ws.Range("K107").Value = 1224.75
$ws.Range("L107").Value = 1497.6666
$ws.Range("M107").Value = 695.25
$ws.Range("N107").Value = -5337.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3478.4
$ws.Range("I31").Value = 2699.5
$ws.Range("J31").Value = 3997.6667
$ws.Range("K31").Value = 2699.5
$ws.Range("L31").Value = 3997.6667
$ws.Range("M31").Value = -2404.5
$ws.Range("N31").Value = -4587.6667
$ws.Range("H34").Value = 3478.4
$ws.Range("I34").Value = 2699.5
$ws.Range("J34").Value = 3997.6667
$ws.Range("K34").Value = 2699.5
$ws.Range("L34").Value = 3997.6667
$ws.Range("M34").Value = -2497.5
$ws.Range("N34").Value = -4401.6667
$ws.Range("H58").Value = 2336.7144
$ws.Range("J58").Value = 2854.3333
$ws.Range("L58").Value = 2854.3333
$ws.Range("N58").Value = -3260.3333
$ws.Range("H136").Value = 2336.7144
$ws.Range("J136").Value = 2854.3333
$ws.Range("L136").Value = 8562.999899999999
$ws.Range("N136").Value = -13662.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 350
$ws.Range("I60").Value = 350
$ws.Range("K60").Value = 1050
$ws.Range("M60").Value = -799
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H118").Value = 1970.5883
$ws.Range("I118").Value = 1970.5883
$ws.Range("K118").Value = 5911.7649
$ws.Range("M118").Value = -4668.7649
$ws.Range("H122").Value = 894
$ws.Range("I122").Value = 890
$ws.Range("J122").Value = 898
$ws.Range("K122").Value = 8010
$ws.Range("L122").Value = 8082
$ws.Range("M122").Value = -5560
$ws.Range("N122").Value = -12982

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1050485.2
$ws.Range("I3").Value = 917149.8
$ws.Range("J3").Value = 1250488.5
$ws.Range("K3").Value = 917149.8
$ws.Range("L3").Value = 1250488.5
$ws.Range("M3").Value = -917033.8
$ws.Range("N3").Value = -1250720.5
$ws.Range("H7").Value = 14038333
$ws.Range("I7").Value = 19381666
$ws.Range("J7").Value = 3351666.8
$ws.Range("K7").Value = 19381666
$ws.Range("L7").Value = 3351666.8
$ws.Range("M7").Value = -19381554
$ws.Range("N7").Value = -3351890.8
$ws.Range("H8").Value = 14038333
$ws.Range("I8").Value = 19381666
$ws.Range("J8").Value = 3351666.8
$ws.Range("K8").Value = 19381666
$ws.Range("L8").Value = 3351666.8
$ws.Range("M8").Value = -19381527
$ws.Range("N8").Value = -3351944.8
$ws.Range("H10").Value = 13668933
$ws.Range("I10").Value = 20500000
$ws.Range("J10").Value = 6800
$ws.Range("K10").Value = 20500000
$ws.Range("L10").Value = 6800
$ws.Range("M10").Value = -20499831
$ws.Range("N10").Value = -7138
$ws.Range("H11").Value = 1787625.4
$ws.Range("I11").Value = 2454200.5
$ws.Range("K11").Value = 2454200.5
$ws.Range("M11").Value = -2454061.5
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H69").Value = 49999
$ws.Range("J69").Value = 49999
$ws.Range("L69").Value = 49999
$ws.Range("N69").Value = -51497
$ws.Range("H72").Value = 49999
$ws.Range("J72").Value = 49999
$ws.Range("L72").Value = 149997
$ws.Range("N72").Value = -157485
$ws.Range("H102").Value = 2122.875
$ws.Range("I102").Value = 2996
$ws.Range("K102").Value = 2996
$ws.Range("M102").Value = -1374
$ws.Range("H122").Value = 2695.158
$ws.Range("I122").Value = 2856.75
$ws.Range("K122").Value = 8570.25
$ws.Range("M122").Value = -6120.25
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 7108.76
$ws.Range("I132").Value = 7636.45
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 22909.35
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -20379.35
$ws.Range("N132").Value = -20054
$ws.Range("H134").Value = 36220.832
$ws.Range("J134").Value = 36220.832
$ws.Range("L134").Value = 108662.496
$ws.Range("N134").Value = -113732.496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10128
$ws.Range("I7").Value = 13333.667
$ws.Range("K7").Value = 13333.667
$ws.Range("M7").Value = -13221.667
$ws.Range("H40").Value = 4416.0835
$ws.Range("I40").Value = 3888.111
$ws.Range("K40").Value = 3888.111
$ws.Range("M40").Value = -3752.111
$ws.Range("H110").Value = 34000
$ws.Range("J110").Value = 34000
$ws.Range("L110").Value = 34000
$ws.Range("N110").Value = -42180
$ws.Range("H126").Value = 10128
$ws.Range("I126").Value = 13333.667
$ws.Range("K126").Value = 40001.001
$ws.Range("M126").Value = -37531.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 3003
$ws.Range("I9").Value = 4006
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 4006
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = -3866
$ws.Range("N9").Value = -2280
$ws.Range("H12").Value = 15000
$ws.Range("J12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("N12").Value = -15284
$ws.Range("H126").Value = 2632.889
$ws.Range("I126").Value = 2185.1428
$ws.Range("K126").Value = 6555.428400000001
$ws.Range("M126").Value = -4085.428400000001
$ws.Range("H132").Value = 5650
$ws.Range("I132").Value = 5650
$ws.Range("K132").Value = 16950
$ws.Range("M132").Value = -14420
